$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header: "Prompts" -> "Queries"
$ws.Range("D1").Value = "Queries"

# Fix typo in context text: "football," -> "football."
$ws.Range("C3").Value = "John is in the play ground.`nJohn picked up the football.`nBob went to the kitchen."

# Update the active selection to D6
$ws.Range("D6").Select()
